# The "Πωλήσεις Έκπτωση 1" summary figures in row 7 (SalesQuantity / Turnover
# totals, cells styled with xf "5") were shifted one column to the right:
#   J7 (4)    -> K7
#   K7 (8.59) -> L7
# Move K7 first so it lands in the now-vacated L7, then move J7 into K7.
# Using Range.Cut (with an explicit destination) on single-cell ranges carries
# both the value and the existing cell style along with it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K7").Cut($ws.Range("L7"))
$ws.Range("J7").Cut($ws.Range("K7"))
$ws.Range("J7").Clear()

# Last-saved cursor/selection moved to J12.
[void]$ws.Range("J12").Select()
